# Update numeric values in result_data_KNN.xlsx ("Update Name of Algo")
# Applies the specific cell value corrections identified in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -12.377
$ws.Range("A8").Value = -22.071
$ws.Range("A10").Value = -21.616
$ws.Range("C11").Value = -12.868
$ws.Range("A12").Value = -21.58
$ws.Range("C12").Value = -10.76
$ws.Range("C15").Value = -13.364
$ws.Range("C17").Value = -13.36
$ws.Range("A18").Value = -21.886
$ws.Range("A25").Value = -21.816
$ws.Range("C26").Value = -13.173
$ws.Range("C27").Value = -13.693
$ws.Range("C28").Value = -13.191
$ws.Range("C32").Value = -12.972
$ws.Range("A37").Value = -19.775
$ws.Range("C37").Value = -12.638
$ws.Range("C41").Value = -12.375
$ws.Range("C47").Value = -12.881
$ws.Range("C51").Value = -11.371
$ws.Range("A55").Value = -21.868
$ws.Range("C65").Value = -12.1
$ws.Range("A68").Value = -21.562
$ws.Range("C73").Value = -12.466
$ws.Range("A77").Value = -20.71
$ws.Range("A78").Value = -20.118
$ws.Range("A79").Value = -21.143
$ws.Range("A80").Value = -20.239
$ws.Range("A81").Value = -21.818
$ws.Range("A82").Value = -22.105
$ws.Range("A84").Value = -22.073
$ws.Range("C84").Value = -13.924
$ws.Range("C85").Value = -12.132
$ws.Range("C89").Value = -13.617
$ws.Range("C93").Value = -10.901
$ws.Range("C95").Value = -11.682
$ws.Range("C98").Value = -13.329
$ws.Range("C99").Value = -11.64
$ws.Range("A101").Value = -20.326
$ws.Range("C101").Value = -12.514
$ws.Range("A102").Value = -20.104
$ws.Range("C102").Value = -12.287
